# Updates the cryptos price list (Price / Volume(1h) columns, and two
# name-swap corrections) to match the latest scrape, per the authoritative
# diff of xl/worksheets/sheet1.xml.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All values in this sheet (Coin, Link, Price, Volume) are stored as plain
# text, even when the text looks like a plain number (e.g. "0.658"). Force
# the cell to Text format before writing such values so Excel does not
# silently reinterpret them as a numeric value (which would change "0.658"
# into a binary double and lose the original text formatting).
function Set-TextValue($range, [string]$value) {
    if ($value -match '^\s*-?\d+(\.\d+)?\s*$') {
        $range.NumberFormat = "@"
    }
    $range.Value = $value
}

# Row 2
Set-TextValue $ws.Range("D2") "36.453.23"
Set-TextValue $ws.Range("E2") "  -1.11%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.034.93"
Set-TextValue $ws.Range("E3") "  -0.78%  "

# Row 4
Set-TextValue $ws.Range("E4") "  +0.03%  "

# Row 5
Set-TextValue $ws.Range("D5") "243.99"
Set-TextValue $ws.Range("E5") "  -0.47%  "

# Row 6
Set-TextValue $ws.Range("D6") "0.658"
Set-TextValue $ws.Range("E6") "  +0.80%  "

# Row 7
Set-TextValue $ws.Range("E7") "  +0.01%  "

# Row 8
Set-TextValue $ws.Range("D8") "53.31"
Set-TextValue $ws.Range("E8") "  -6.71%  "

# Row 9
Set-TextValue $ws.Range("D9") "62.46"
Set-TextValue $ws.Range("E9") "  +5.66%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.360"
Set-TextValue $ws.Range("E10") "  -2.11%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.0740"
Set-TextValue $ws.Range("E11") "  -4.51%  "

# Row 12
Set-TextValue $ws.Range("E12") "  -3.94%  "

# Row 13
Set-TextValue $ws.Range("D13") "0.926"
Set-TextValue $ws.Range("E13") "  +6.39%  "

# Row 14
Set-TextValue $ws.Range("D14") "14.30"
Set-TextValue $ws.Range("E14") "  -5.22%  "

# Row 15
Set-TextValue $ws.Range("D15") "2.333.61"
Set-TextValue $ws.Range("E15") "  -0.77%  "

# Row 16
Set-TextValue $ws.Range("D16") "5.35"
Set-TextValue $ws.Range("E16") "  -3.77%  "

# Row 17
Set-TextValue $ws.Range("D17") "2.107.34"
Set-TextValue $ws.Range("E17") "  +3.98%  "

# Row 18
Set-TextValue $ws.Range("D18") "36.344.43"
Set-TextValue $ws.Range("E18") "  -1.30%  "

# Row 19
Set-TextValue $ws.Range("D19") "16.87"
Set-TextValue $ws.Range("E19") "  -2.49%  "

# Row 20
Set-TextValue $ws.Range("D20") "70.84"
Set-TextValue $ws.Range("E20") "  -2.90%  "

# Row 21
Set-TextValue $ws.Range("D21") "0.0₃0847"
Set-TextValue $ws.Range("E21") "  -4.41%  "

# Row 22
Set-TextValue $ws.Range("D22") "236.69"
Set-TextValue $ws.Range("E22") "  +0.46%  "

# Row 23
Set-TextValue $ws.Range("D23") "5.12"
Set-TextValue $ws.Range("E23") "  -4.81%  "

# Row 24
Set-TextValue $ws.Range("E24") "  +0.06%  "

# Row 25
Set-TextValue $ws.Range("E25") "  -3.96%  "

# Row 26
Set-TextValue $ws.Range("D26") "2.21"
Set-TextValue $ws.Range("E26") "  +0.53%  "

# Row 27
Set-TextValue $ws.Range("B27") "Cosmos"
Set-TextValue $ws.Range("C27") "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D27") "9.11"
Set-TextValue $ws.Range("E27") "  -9.72%  "

# Row 28
Set-TextValue $ws.Range("B28") "Monero"
Set-TextValue $ws.Range("C28") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D28") "163.25"
Set-TextValue $ws.Range("E28") "  -3.12%  "

# Row 29
Set-TextValue $ws.Range("D29") "19.76"
Set-TextValue $ws.Range("E29") "  -0.88%  "

# Row 30
Set-TextValue $ws.Range("E30") "  -3.09%  "

# Row 31
Set-TextValue $ws.Range("E31") "  +4.07%  "

# Row 32
Set-TextValue $ws.Range("D32") "4.93"
Set-TextValue $ws.Range("E32") "  -9.54%  "

# Row 33
Set-TextValue $ws.Range("D33") "0.0592"
Set-TextValue $ws.Range("E33") "  -3.03%  "

# Row 34
Set-TextValue $ws.Range("D34") "4.38"
Set-TextValue $ws.Range("E34") "  -9.06%  "

# Row 35
Set-TextValue $ws.Range("E35") "  +0.11%  "

# Row 36
Set-TextValue $ws.Range("D36") "0.0861"
Set-TextValue $ws.Range("E36") "  +0.42%  "

# Row 37
Set-TextValue $ws.Range("E37") "  -1.43%  "

# Row 38
Set-TextValue $ws.Range("D38") "2.18"
Set-TextValue $ws.Range("E38") "  -5.99%  "

# Row 39
Set-TextValue $ws.Range("D39") "1.22"
Set-TextValue $ws.Range("E39") "  -6.46%  "

# Row 40
Set-TextValue $ws.Range("D40") "4.87"
Set-TextValue $ws.Range("E40") "  +0.04%  "

# Row 41
Set-TextValue $ws.Range("E41") "  -3.66%  "

# Row 42
Set-TextValue $ws.Range("E42") "  -4.48%  "

# Row 43
Set-TextValue $ws.Range("E43") "  -4.30%  "

# Row 44
Set-TextValue $ws.Range("D44") "92.65"
Set-TextValue $ws.Range("E44") "  -4.10%  "

# Row 45
Set-TextValue $ws.Range("D45") "0.0891"
Set-TextValue $ws.Range("E45") "  -6.03%  "

# Row 46
Set-TextValue $ws.Range("D46") "1.368.36"
Set-TextValue $ws.Range("E46") "  +5.00%  "

# Row 47
Set-TextValue $ws.Range("B47") "FraxShare"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D47") "7.44"
Set-TextValue $ws.Range("E47") "  +10.56%  "

# Row 48
Set-TextValue $ws.Range("B48") "InjectiveProtocol"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D48") "15.59"
Set-TextValue $ws.Range("E48") "  -5.95%  "

# Row 49
Set-TextValue $ws.Range("E49") "  +2.38%  "

# Row 50
Set-TextValue $ws.Range("D50") "2.24"
Set-TextValue $ws.Range("E50") "  -4.68%  "

# Row 51
Set-TextValue $ws.Range("D51") "2.222.23"

